$wb = $excel.ActiveWorkbook

# A never-touched far-away cell used purely as a "blank / no style" style
# donor (see the PasteSpecial-formats trick below).
$blankDonorSheet = $wb.Worksheets.Item("2021-Q3")
$blank = $blankDonorSheet.Range("Z999")

# Helper-less pattern repeated by hand below (the COM bridge has no script
# functions in this environment, so everything is inlined):
#   1. $cell.Value = <text>              -> writes the literal text
#   2. $blank.Copy(); $cell.PasteSpecial(-4122)
#                                         -> clears any number formatting /
#                                            quote-prefix style picked up
#                                            while typing a numeric-looking
#                                            string, WITHOUT touching the
#                                            value that was just written.

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q3" (and before "总计")
# ---------------------------------------------------------------------------
$ws2021 = $wb.Worksheets.Item("2021-Q3")

$ws2022 = $wb.Worksheets.Add([Type]::Missing, $ws2021)
$ws2022.Name = "2022-Q1"

# Re-fetch "总计" AFTER the insert: sheet references resolve by tab position,
# so a handle captured before the insert would now point at the new sheet.
$wsTotal = $wb.Worksheets.Item("总计")

# Clone the header / row-label formatting used on "总计" (style used for
# B1:D1 and A2) onto the new sheet - that is the style the target file uses
# for this sheet's header row and its "A" index column.
$wsTotal.Range("B1:D1").Copy($ws2022.Range("B1:D1"))
$wsTotal.Range("B1:D1").Copy($ws2022.Range("E1:G1"))
$wsTotal.Range("B1").Copy($ws2022.Range("H1"))
$wsTotal.Range("A2").Copy($ws2022.Range("A2"))
$wsTotal.Range("A2").Copy($ws2022.Range("A3"))

# Header row text
$ws2022.Cells.Item(1, 2).Value = "基金代码"
$ws2022.Cells.Item(1, 3).Value = "基金名称"
$ws2022.Cells.Item(1, 4).Value = "基金规模"
$ws2022.Cells.Item(1, 5).Value = "股票总仓位"
$ws2022.Cells.Item(1, 6).Value = "仓位占比"
$ws2022.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws2022.Cells.Item(1, 8).Value = "仓位排名"

# Row 2 - fund 004634
$ws2022.Cells.Item(2, 1).Value = 0
$ws2022.Cells.Item(2, 2).Value = "'004634"
$blank.Copy()
$ws2022.Cells.Item(2, 2).PasteSpecial(-4122)
$ws2022.Cells.Item(2, 3).Value = "新疆前海联合泳涛灵活配置混合A"
$ws2022.Cells.Item(2, 4).Value = "'1.33"
$blank.Copy()
$ws2022.Cells.Item(2, 4).PasteSpecial(-4122)
$ws2022.Cells.Item(2, 5).Value = "'89.65"
$blank.Copy()
$ws2022.Cells.Item(2, 5).PasteSpecial(-4122)
$ws2022.Cells.Item(2, 6).Value = "'4.13"
$blank.Copy()
$ws2022.Cells.Item(2, 6).PasteSpecial(-4122)
$ws2022.Cells.Item(2, 7).Value = "'0.0549"
$blank.Copy()
$ws2022.Cells.Item(2, 7).PasteSpecial(-4122)
$ws2022.Cells.Item(2, 8).Value = 10

# Row 3 - fund 007041
$ws2022.Cells.Item(3, 1).Value = 1
$ws2022.Cells.Item(3, 2).Value = "'007041"
$blank.Copy()
$ws2022.Cells.Item(3, 2).PasteSpecial(-4122)
$ws2022.Cells.Item(3, 3).Value = "新疆前海联合泳涛灵活配置混合C"
$ws2022.Cells.Item(3, 4).Value = "'0.00"
$blank.Copy()
$ws2022.Cells.Item(3, 4).PasteSpecial(-4122)
$ws2022.Cells.Item(3, 5).Value = "'89.65"
$blank.Copy()
$ws2022.Cells.Item(3, 5).PasteSpecial(-4122)
$ws2022.Cells.Item(3, 6).Value = "'4.13"
$blank.Copy()
$ws2022.Cells.Item(3, 6).PasteSpecial(-4122)
$ws2022.Cells.Item(3, 7).Value = 0
$ws2022.Cells.Item(3, 8).Value = 10

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet (old row 2 becomes
#    row 3, keeping its original 2021-Q3 values/style).
# ---------------------------------------------------------------------------
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 0.05

$wsTotal.Cells.Item(3, 1).Value = 1

# Restore the original active sheet/tab (the source workbook had "2021-Q3"
# selected and nothing in the diff indicates that should change).
$ws2021.Activate()
